# Update res_bus/vm_pu.xlsx: voltage-magnitude results for the 380 kV case
# (commit: "case with 380 kV done"). Rewrites the B:N value columns for rows 2-25
# and clears the stray H2 cell (column H only ever had data in the old row 2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.172898460946372
$ws.Range("D2").Value = 1.179190980263645
$ws.Range("E2").Value = 1.164525215978536
$ws.Range("F2").Value = 1.175994904095898
$ws.Range("G2").Value = 1
$ws.Range("I2").Value = 1.036733453855621
$ws.Range("J2").Value = 1.1773599188047
$ws.Range("K2").Value = 1.181629437015791
$ws.Range("L2").Value = 1.166995871270124
$ws.Range("M2").Value = 1.178440304681063
$ws.Range("N2").Value = 1.179031904934473

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.17803832402733
$ws.Range("D3").Value = 1.184225271503157
$ws.Range("E3").Value = 1.169149948550072
$ws.Range("F3").Value = 1.180869966940986
$ws.Range("G3").Value = 1
$ws.Range("I3").Value = 1.037192878988106
$ws.Range("J3").Value = 1.182176470244498
$ws.Range("K3").Value = 1.186490328027073
$ws.Range("L3").Value = 1.171445526861641
$ws.Range("M3").Value = 1.183141744957338
$ws.Range("N3").Value = 1.183855296429781

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.181314603963118
$ws.Range("D4").Value = 1.187432952592099
$ws.Range("E4").Value = 1.172096510424995
$ws.Range("F4").Value = 1.183973254604411
$ws.Range("G4").Value = 1
$ws.Range("I4").Value = 1.037478022231475
$ws.Range("J4").Value = 1.185244342682245
$ws.Range("K4").Value = 1.189585965413035
$ws.Range("L4").Value = 1.174278900821462
$ws.Range("M4").Value = 1.186132824099284
$ws.Range("N4").Value = 1.18692752559828

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.182680516165131
$ws.Range("D5").Value = 1.188769951209466
$ws.Range("E5").Value = 1.173324632447065
$ws.Range("F5").Value = 1.185266035466533
$ws.Range("G5").Value = 1
$ws.Range("I5").Value = 1.037595040876778
$ws.Range("J5").Value = 1.186522811396982
$ws.Range("K5").Value = 1.190875889781271
$ws.Range("L5").Value = 1.175459454050475
$ws.Range("M5").Value = 1.187378456484886
$ws.Range("N5").Value = 1.188207809885234

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.182909200304763
$ws.Range("D6").Value = 1.188993776084781
$ws.Range("E6").Value = 1.173530228158815
$ws.Range("F6").Value = 1.185482416454744
$ws.Range("G6").Value = 1
$ws.Range("I6").Value = 1.03761452276508
$ws.Range("J6").Value = 1.186736822860723
$ws.Range("K6").Value = 1.191091812127773
$ws.Range("L6").Value = 1.175657063423414
$ws.Range("M6").Value = 1.187586922193787
$ws.Range("N6").Value = 1.188422125269801

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.181332899731013
$ws.Range("D7").Value = 1.187450862316241
$ws.Range("E7").Value = 1.172112961850427
$ws.Range("F7").Value = 1.183990574794536
$ws.Range("G7").Value = 1
$ws.Range("I7").Value = 1.037479596996601
$ws.Range("J7").Value = 1.185261469379502
$ws.Range("K7").Value = 1.189603246026593
$ws.Range("L7").Value = 1.174294716580449
$ws.Range("M7").Value = 1.186149514206297
$ws.Range("N7").Value = 1.186944676617411

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.174646026271074
$ws.Range("D8").Value = 1.180902922818698
$ws.Range("E8").Value = 1.166097919299112
$ws.Range("F8").Value = 1.177653313839638
$ws.Range("G8").Value = 1
$ws.Range("I8").Value = 1.036891262261581
$ws.Range("J8").Value = 1.178998042341821
$ws.Range("K8").Value = 1.183282740021534
$ws.Range("L8").Value = 1.168509381052595
$ws.Range("M8").Value = 1.180040004120601
$ws.Range("N8").Value = 1.180672354794912

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.162463133237349
$ws.Range("D9").Value = 1.168962930677166
$ws.Range("E9").Value = 1.155128382232715
$ws.Range("F9").Value = 1.166074345691024
$ws.Range("G9").Value = 1
$ws.Range("I9").Value = 1.035759222261601
$ws.Range("J9").Value = 1.167568345648049
$ws.Range("K9").Value = 1.171745139686581
$ws.Range("L9").Value = 1.157945756727379
$ws.Range("M9").Value = 1.168863822809237
$ws.Range("N9").Value = 1.16922642662083

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.154042498111525
$ws.Range("D10").Value = 1.16070334499222
$ws.Range("E10").Value = 1.147539277217807
$ws.Range("F10").Value = 1.158048874015242
$ws.Range("G10").Value = 1
$ws.Range("I10").Value = 1.034936893787768
$ws.Range("J10").Value = 1.159655931831867
$ws.Range("K10").Value = 1.163755518420129
$ws.Range("L10").Value = 1.150628590143588
$ws.Range("M10").Value = 1.161108462723247
$ws.Range("N10").Value = 1.161302776269463

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.150318236932015
$ws.Range("D11").Value = 1.157048706012831
$ws.Range("E11").Value = 1.144181113414466
$ws.Range("F11").Value = 1.154494119900936
$ws.Range("G11").Value = 1
$ws.Range("I11").Value = 1.034563913181882
$ws.Range("J11").Value = 1.156153506130681
$ws.Range("K11").Value = 1.160218326338071
$ws.Range("L11").Value = 1.147388628095247
$ws.Range("M11").Value = 1.157671190239914
$ws.Range("N11").Value = 1.157795376722048

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.148922494883345
$ws.Range("D12").Value = 1.155678816386903
$ws.Range("E12").Value = 1.142922324671092
$ws.Range("F12").Value = 1.153161114461964
$ws.Range("G12").Value = 1
$ws.Range("I12").Value = 1.03442275217299
$ws.Range("J12").Value = 1.154840459528738
$ws.Range("K12").Value = 1.158892157622039
$ws.Range("L12").Value = 1.146173824592848
$ws.Range("M12").Value = 1.156381914121668
$ws.Range("N12").Value = 1.156480465443321

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.149222457769985
$ws.Range("D13").Value = 1.15597323407491
$ws.Range("E13").Value = 1.14319286576071
$ws.Range("F13").Value = 1.153447630108703
$ws.Range("G13").Value = 1
$ws.Range("I13").Value = 1.034453151555342
$ws.Range("J13").Value = 1.155122670076862
$ws.Range("K13").Value = 1.159177192465995
$ws.Range("L13").Value = 1.14643492688275
$ws.Range("M13").Value = 1.156659045325715
$ws.Range("N13").Value = 1.156763076762794

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.150203121065715
$ws.Range("D14").Value = 1.156935727234604
$ws.Range("E14").Value = 1.14407729802814
$ws.Range("F14").Value = 1.154384194382054
$ws.Range("G14").Value = 1
$ws.Range("I14").Value = 1.034552298681632
$ws.Range("J14").Value = 1.156045219682435
$ws.Range("K14").Value = 1.160108959508403
$ws.Range("L14").Value = 1.147288446841046
$ws.Range("M14").Value = 1.157564877457266
$ws.Range("N14").Value = 1.157686936494625

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.150805678530388
$ws.Range("D15").Value = 1.157527088638154
$ws.Range("E15").Value = 1.144620694601321
$ws.Range("F15").Value = 1.154959551507918
$ws.Range("G15").Value = 1
$ws.Range("I15").Value = 1.034613036938153
$ws.Range("J15").Value = 1.156612011507045
$ws.Range("K15").Value = 1.160681402653132
$ws.Range("L15").Value = 1.147812808276191
$ws.Range("M15").Value = 1.158121312013103
$ws.Range("N15").Value = 1.158254533228636

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.15428796512109
$ws.Range("D16").Value = 1.160944189746538
$ws.Range("E16").Value = 1.147760579919141
$ws.Range("F16").Value = 1.158283058900293
$ws.Range("G16").Value = 1
$ws.Range("I16").Value = 1.034961284422651
$ws.Range("J16").Value = 1.15988671617739
$ws.Range("K16").Value = 1.163988581317732
$ws.Range("L16").Value = 1.150842058698251
$ws.Range("M16").Value = 1.161334862824321
$ws.Range("N16").Value = 1.161533888355228

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.156450959343791
$ws.Range("D17").Value = 1.163066270466847
$ws.Range("E17").Value = 1.149710453699905
$ws.Range("F17").Value = 1.160346036564162
$ws.Range("G17").Value = 1
$ws.Range("I17").Value = 1.035175150061517
$ws.Range("J17").Value = 1.161919994576496
$ws.Range("K17").Value = 1.16604186705219
$ws.Range("L17").Value = 1.152722663127862
$ws.Range("M17").Value = 1.163329016234961
$ws.Range("N17").Value = 1.163570054242881

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.157705111083362
$ws.Range("D18").Value = 1.1642965465138
$ws.Range("E18").Value = 1.150840874274837
$ws.Range("F18").Value = 1.161541694801962
$ws.Range("G18").Value = 1
$ws.Range("I18").Value = 1.035298268632654
$ws.Range("J18").Value = 1.163098654179333
$ws.Range("K18").Value = 1.167232068333837
$ws.Range("L18").Value = 1.153812722683647
$ws.Range("M18").Value = 1.164484579848642
$ws.Range("N18").Value = 1.164750387677548

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.158131493880987
$ws.Range("D19").Value = 1.164714786105393
$ws.Range("E19").Value = 1.151225164254658
$ws.Range("F19").Value = 1.161948106011007
$ws.Range("G19").Value = 1
$ws.Range("I19").Value = 1.035339975468836
$ws.Range("J19").Value = 1.163499323889588
$ws.Range("K19").Value = 1.167636651929979
$ws.Range("L19").Value = 1.154183257471127
$ws.Range("M19").Value = 1.164877327821977
$ws.Range("N19").Value = 1.165151626384749

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.156219669477513
$ws.Range("D20").Value = 1.162839371364758
$ws.Range("E20").Value = 1.14950196938602
$ws.Range("F20").Value = 1.160125493511911
$ws.Range("G20").Value = 1
$ws.Range("I20").Value = 1.035152373006177
$ws.Range("J20").Value = 1.161702604311301
$ws.Range("K20").Value = 1.16582234351034
$ws.Range("L20").Value = 1.152521606265762
$ws.Range("M20").Value = 1.163115852172341
$ws.Range("N20").Value = 1.163352355258574

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.149914687554645
$ws.Range("D21").Value = 1.156652644534481
$ws.Range("E21").Value = 1.143817174912362
$ws.Range("F21").Value = 1.154108753046734
$ws.Range("G21").Value = 1
$ws.Range("I21").Value = 1.034523175296729
$ws.Range("J21").Value = 1.15577389079014
$ws.Range("K21").Value = 1.159834922126629
$ws.Range("L21").Value = 1.147037424324511
$ws.Range("M21").Value = 1.157298483339473
$ws.Range("N21").Value = 1.15741522228418

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.145878479476662
$ws.Range("D22").Value = 1.152690744411885
$ws.Range("E22").Value = 1.140176539558095
$ws.Range("F22").Value = 1.150252474914388
$ws.Range("G22").Value = 1
$ws.Range("I22").Value = 1.034112375375829
$ws.Range("J22").Value = 1.151975985444667
$ws.Range("K22").Value = 1.15599889524057
$ws.Range("L22").Value = 1.14352339054932
$ws.Range("M22").Value = 1.15356809399958
$ws.Range("N22").Value = 1.153611923477491

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.148025206350399
$ws.Range("D23").Value = 1.154798080350242
$ws.Range("E23").Value = 1.142113010152756
$ws.Range("F23").Value = 1.152303935519696
$ws.Range("G23").Value = 1
$ws.Range("I23").Value = 1.034331616731898
$ws.Range("J23").Value = 1.153996209377484
$ws.Range("K23").Value = 1.158039445506821
$ws.Range("L23").Value = 1.14539269859393
$ws.Range("M23").Value = 1.155552763204379
$ws.Range("N23").Value = 1.155635016360015

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.156324202499726
$ws.Range("D24").Value = 1.162941920425893
$ws.Range("E24").Value = 1.149596195777754
$ws.Range("F24").Value = 1.160225170993084
$ws.Range("G24").Value = 1
$ws.Range("I24").Value = 1.035162669994981
$ws.Range("J24").Value = 1.161800856182186
$ws.Range("K24").Value = 1.16592155971292
$ws.Range("L24").Value = 1.152612476380862
$ws.Range("M24").Value = 1.163212195253968
$ws.Range("N24").Value = 1.163450746658386

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.165663049625547
$ws.Range("D25").Value = 1.172100238193052
$ws.Range("E25").Value = 1.158010843525314
$ws.Range("F25").Value = 1.169119500237839
$ws.Range("G25").Value = 1
$ws.Range("I25").Value = 1.036063528351127
$ws.Range("J25").Value = 1.170572573868437
$ws.Range("K25").Value = 1.174778165847196
$ws.Range("L25").Value = 1.160723088580781
$ws.Range("M25").Value = 1.171804613901682
$ws.Range("N25").Value = 1.172234921189882

# H2 no longer holds data in the updated results (only B2:N2 except H2 are populated)
$ws.Range("H2").ClearContents()

Write-Output "applied vm_pu updates for Case_0_173 (380 kV)"
